$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Name" header in B1, matching the formatting already used for A1 ("Roll Number")
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B1").Value = "Name"

# New data row: roll number (kept as text, like the existing roll-number column) + student name
$ws.Range("A2").Value = "'3876"
$ws.Range("B2").Value = "Meharjot"
